$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B8").Value = "mux"
$ws.Range("C8").Value = 0.989
$ws.Range("D8").Value = 1.958
$ws.Range("E8").Value = 3.032
$ws.Range("F8").Value = 4.156
$ws.Range("G8").Value = 5.423
